# Add a new bulleted list item ("abcd") to the end of the "Source Code
# Files" list, right after "SMSConfiguration.aspx" - i.e. add one more
# row to the list (commit message: "Added 1 row").

$d = $word.ActiveDocument

# Move to the very end of the document's main story and start a new
# paragraph there (same as placing the cursor at the end of the last
# line and pressing Enter in Word).
$endRange = $d.Content
$endRange.Collapse(0)          # wdCollapseEnd
$endRange.InsertAfter("`r")
$endRange.Collapse(0)

# The freshly created paragraph inherits the preceding paragraph's
# style/numbering (ListParagraph, same numId/ilvl), exactly like typing
# Enter after the previous list item would in Word. Just add the text.
$endRange.InsertAfter("abcd")

Write-Output "Inserted new list paragraph 'abcd'."
